# Filter the "Data" sheet's table (A1:T58) on the "Cluster" column (column M,
# the 13th column) to show only rows whose Cluster is "Hardware und Logik" or
# "Prozessoren und Architekturen" (commit: "working on 5-1 Trends").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A1:T58").AutoFilter(13, @("Hardware und Logik", "Prozessoren und Architekturen"), 7)
